$p = $ppt.ActivePresentation

# Insert a new blank slide as slide 2 (after the existing slide 1)
$s = $p.Slides.Add(2, 12)  # 12 = ppLayoutBlank
